$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing rows 2-27: new A (id) values
$ws.Range("A2").Value = 10001303
$ws.Range("A3").Value = 10002303
$ws.Range("A4").Value = 10003303
$ws.Range("A5").Value = 10004303
$ws.Range("A6").Value = 10005303
$ws.Range("A7").Value = 10006303
$ws.Range("A8").Value = 10007303
$ws.Range("A9").Value = 10008303
$ws.Range("A10").Value = 10009303
$ws.Range("A11").Value = 10010303
$ws.Range("A12").Value = 10011303
$ws.Range("A13").Value = 10012303
$ws.Range("A14").Value = 10013303
$ws.Range("A15").Value = 10014303
$ws.Range("A16").Value = 10015303
$ws.Range("A17").Value = 10016303
$ws.Range("A18").Value = 10017303
$ws.Range("A19").Value = 10018303
$ws.Range("A20").Value = 10019303
$ws.Range("A21").Value = 10020303
$ws.Range("A22").Value = 10021303
$ws.Range("A23").Value = 10022303
$ws.Range("A24").Value = 10023303
$ws.Range("A25").Value = 10024303
$ws.Range("A26").Value = 10025303
$ws.Range("A27").Value = 10026303

# Add new rows 28-53 (B, C, D same pattern as rows 2-27, new A ids, D=404)
$ws.Range("A28").Value = 10027404
$ws.Range("B28").Value = "AKL"
$ws.Range("C28").Value = "X AKL 1"
$ws.Range("D28").Value = 404
$ws.Range("A29").Value = 10028404
$ws.Range("B29").Value = "AKL"
$ws.Range("C29").Value = "X AKL 2"
$ws.Range("D29").Value = 404
$ws.Range("A30").Value = 10029404
$ws.Range("B30").Value = "PM"
$ws.Range("C30").Value = "X PM 1"
$ws.Range("D30").Value = 404
$ws.Range("A31").Value = 10030404
$ws.Range("B31").Value = "MPLB"
$ws.Range("C31").Value = "X MPLB 1"
$ws.Range("D31").Value = 404
$ws.Range("A32").Value = 10031404
$ws.Range("B32").Value = "MPLB"
$ws.Range("C32").Value = "X MPLB 2"
$ws.Range("D32").Value = 404
$ws.Range("A33").Value = 10032404
$ws.Range("B33").Value = "TJKT"
$ws.Range("C33").Value = "X TJKT 1"
$ws.Range("D33").Value = 404
$ws.Range("A34").Value = 10033404
$ws.Range("B34").Value = "TJKT"
$ws.Range("C34").Value = "X TJKT 2"
$ws.Range("D34").Value = 404
$ws.Range("A35").Value = 10034404
$ws.Range("B35").Value = "TJKT"
$ws.Range("C35").Value = "X TJKT 3"
$ws.Range("D35").Value = 404
$ws.Range("A36").Value = 10035404
$ws.Range("B36").Value = "DKV"
$ws.Range("C36").Value = "X DKV 1"
$ws.Range("D36").Value = 404
$ws.Range("A37").Value = 10036404
$ws.Range("B37").Value = "AKL"
$ws.Range("C37").Value = "XI AKL 1"
$ws.Range("D37").Value = 404
$ws.Range("A38").Value = 10037404
$ws.Range("B38").Value = "AKL"
$ws.Range("C38").Value = "XI AKL 2"
$ws.Range("D38").Value = 404
$ws.Range("A39").Value = 10038404
$ws.Range("B39").Value = "PM"
$ws.Range("C39").Value = "XI PM 1"
$ws.Range("D39").Value = 404
$ws.Range("A40").Value = 10039404
$ws.Range("B40").Value = "MPLB"
$ws.Range("C40").Value = "XI MPLB 1"
$ws.Range("D40").Value = 404
$ws.Range("A41").Value = 10040404
$ws.Range("B41").Value = "MPLB"
$ws.Range("C41").Value = "XI MPLB 2"
$ws.Range("D41").Value = 404
$ws.Range("A42").Value = 10041404
$ws.Range("B42").Value = "TJKT"
$ws.Range("C42").Value = "XI TJKT 1"
$ws.Range("D42").Value = 404
$ws.Range("A43").Value = 10042404
$ws.Range("B43").Value = "TJKT"
$ws.Range("C43").Value = "XI TJKT 2"
$ws.Range("D43").Value = 404
$ws.Range("A44").Value = 10043404
$ws.Range("B44").Value = "TJKT"
$ws.Range("C44").Value = "XI TJKT 3"
$ws.Range("D44").Value = 404
$ws.Range("A45").Value = 10044404
$ws.Range("B45").Value = "AKL"
$ws.Range("C45").Value = "XII AKL 1"
$ws.Range("D45").Value = 404
$ws.Range("A46").Value = 10045404
$ws.Range("B46").Value = "AKL"
$ws.Range("C46").Value = "XII AKL 2"
$ws.Range("D46").Value = 404
$ws.Range("A47").Value = 10046404
$ws.Range("B47").Value = "PM"
$ws.Range("C47").Value = "XII PM 1"
$ws.Range("D47").Value = 404
$ws.Range("A48").Value = 10047404
$ws.Range("B48").Value = "PM"
$ws.Range("C48").Value = "XII PM 2"
$ws.Range("D48").Value = 404
$ws.Range("A49").Value = 10048404
$ws.Range("B49").Value = "MPLB"
$ws.Range("C49").Value = "XII MPLB 1"
$ws.Range("D49").Value = 404
$ws.Range("A50").Value = 10049404
$ws.Range("B50").Value = "MPLB"
$ws.Range("C50").Value = "XII MPLB 2"
$ws.Range("D50").Value = 404
$ws.Range("A51").Value = 10050404
$ws.Range("B51").Value = "TJKT"
$ws.Range("C51").Value = "XII TJKT 1"
$ws.Range("D51").Value = 404
$ws.Range("A52").Value = 10051404
$ws.Range("B52").Value = "TJKT"
$ws.Range("C52").Value = "XII TJKT 2"
$ws.Range("D52").Value = 404
$ws.Range("A53").Value = 10052404
$ws.Range("B53").Value = "TJKT"
$ws.Range("C53").Value = "XII TJKT 3"
$ws.Range("D53").Value = 404

# Apply General number format to column A data cells (matches new style index with applyNumberFormat)
$ws.Range("A2:A53").NumberFormat = "General"

# Update sheet view
$ws.Application.ActiveWindow.ScrollRow = 14
$ws.Range("D41").Select()
